$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "SNOW-773795"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2/26/2026"
$ws.Range("C8").Value = "desk"
$ws.Range("D8").Value = "d@d.com"
$ws.Range("E8").Value = "d"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "Weekend Getaway"
$ws.Range("H8").Value = 12800
$ws.Range("I8").Value = 12800
$ws.Range("J8").Value = "Confirmed"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "2/23/2026"
$ws.Range("L8").Value = ""
